# fixing build error (24/12/2023)
# Update a handful of text cells on the "Tháng 12" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tháng 12")

$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "456456"
$ws.Range("C7").Value = "Lương tháng"
$ws.Range("D7").Value = "hfgh"

$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "65656756"
$ws.Range("J13").Value = "rthgfhfghfgh"
